# Auto-generated Excel COM-interop script applying market-price value updates
# to the Jenova_Profits leve-profit tracking workbook (scheduled Universalis refresh).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 113 (hunk 0)
$ws.Range("H113").Value = 2235
$ws.Range("I113").Value = 2235
$ws.Range("K113").Value = 2235
$ws.Range("M113").Value = 1019
# Row 133 (hunk 1)
$ws.Range("H133").Value = 49959.098
$ws.Range("J133").Value = 49959.098
$ws.Range("L133").Value = 49959.098
$ws.Range("N133").Value = -60079.098
# Row 135 (hunk 2)
$ws.Range("H135").Value = 1795.9375
$ws.Range("I135").Value = 1520.3572
$ws.Range("K135").Value = 13683.2148
$ws.Range("M135").Value = -11148.2148
# Row 137 (hunk 3)
$ws.Range("H137").Value = 1569.7457
$ws.Range("I137").Value = 1241.4318
$ws.Range("K137").Value = 3724.2954
$ws.Range("M137").Value = -1174.2954
# Row 138 (hunk 4)
$ws.Range("H138").Value = 6590.2334
$ws.Range("J138").Value = 8274.927
$ws.Range("L138").Value = 24824.781
$ws.Range("N138").Value = -35104.781

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (hunk 5)
$ws.Range("H2").Value = 41506.926
$ws.Range("I2").Value = 55513.8
$ws.Range("K2").Value = 55513.8
$ws.Range("M2").Value = -55400.8
# Row 3 (hunk 6)
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").Value = ""
# Row 4 (hunk 7)
$ws.Range("H4").Value = 300
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = ""
# Row 5 (hunk 8)
$ws.Range("H5").Value = 175
$ws.Range("I5").Value = 133.66667
$ws.Range("K5").Value = 133.66667
$ws.Range("M5").Value = -21.66667000000001
# Row 14 (hunk 9)
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = ""
# Row 16 (hunk 10)
$ws.Range("H16").Value = 703
$ws.Range("I16").Value = 703
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 703
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -416
$ws.Range("N16").Value = ""
# Row 45 (hunk 11)
$ws.Range("H45").Value = 1875.3928
$ws.Range("I45").Value = 1544.4736
$ws.Range("J45").Value = 2574
$ws.Range("K45").Value = 1544.4736
$ws.Range("L45").Value = 2574
$ws.Range("M45").Value = -1167.4736
$ws.Range("N45").Value = -3328
# Row 61 (hunk 12)
$ws.Range("H61").Value = 2870.2
$ws.Range("I61").Value = 2771.4
$ws.Range("K61").Value = 2771.4
$ws.Range("M61").Value = -2559.4
# Row 74 (hunk 13)
$ws.Range("H74").Value = 1707.9048
$ws.Range("I74").Value = 1768.1
$ws.Range("K74").Value = 1768.1
$ws.Range("M74").Value = -894.0999999999999
# Row 77 (hunk 14)
$ws.Range("H77").Value = 1707.9048
$ws.Range("I77").Value = 1768.1
$ws.Range("K77").Value = 8840.5
$ws.Range("M77").Value = -4472.5
# Row 116 (hunk 15)
$ws.Range("H116").Value = 41506.926
$ws.Range("I116").Value = 55513.8
$ws.Range("K116").Value = 55513.8
$ws.Range("M116").Value = -53219.8
# Row 123 (hunk 16)
$ws.Range("H123").Value = 68745.25
$ws.Range("I123").Value = 20000
$ws.Range("J123").Value = 84993.664
$ws.Range("K123").Value = 20000
$ws.Range("L123").Value = 84993.664
$ws.Range("N123").Value = -94793.664
$ws.Range("M123").Value = -15100
# Row 132 (hunk 17)
$ws.Range("H132").Value = 4749.4165
$ws.Range("I132").Value = 3654.0908
$ws.Range("J132").Value = 7159.1333
$ws.Range("K132").Value = 10962.2724
$ws.Range("L132").Value = 21477.3999
$ws.Range("M132").Value = -8432.2724
$ws.Range("N132").Value = -26537.3999
# Row 136 (hunk 18)
$ws.Range("H136").Value = 2870.2
$ws.Range("I136").Value = 2771.4
$ws.Range("K136").Value = 8314.200000000001
$ws.Range("M136").Value = -5764.200000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (hunk 19)
$ws.Range("H3").Value = 41506.926
$ws.Range("I3").Value = 55513.8
$ws.Range("K3").Value = 55513.8
$ws.Range("M3").Value = -55399.8
# Row 4 (hunk 20)
$ws.Range("H4").Value = 175
$ws.Range("I4").Value = 133.66667
$ws.Range("K4").Value = 133.66667
$ws.Range("M4").Value = -18.66667000000001
# Row 20 (hunk 21)
$ws.Range("H20").Value = 3487.125
$ws.Range("I20").Value = 3299.9
$ws.Range("J20").Value = 3799.1667
$ws.Range("K20").Value = 3299.9
$ws.Range("L20").Value = 3799.1667
$ws.Range("M20").Value = -3052.9
$ws.Range("N20").Value = -4293.1667
# Row 22 (hunk 22)
$ws.Range("H22").Value = 258.4
$ws.Range("I22").Value = 258.4
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 258.4
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -85.39999999999998
$ws.Range("N22").Value = ""
# Row 134 (hunk 23)
$ws.Range("H134").Value = 23168.96
$ws.Range("I134").Value = 3234.8157
$ws.Range("K134").Value = 9704.447100000001
$ws.Range("M134").Value = -7169.447100000001

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7 (hunk 24)
$ws.Range("H7").Value = 402.72
$ws.Range("I7").Value = 395.13635
$ws.Range("K7").Value = 395.13635
$ws.Range("M7").Value = -282.13635
# Row 16 (hunk 25)
$ws.Range("H16").Value = 3429.0417
$ws.Range("J16").Value = 3719.625
$ws.Range("L16").Value = 3719.625
$ws.Range("N16").Value = -4293.625
# Row 22 (hunk 26)
$ws.Range("H22").Value = 397.83334
$ws.Range("I22").Value = 529
$ws.Range("K22").Value = 529
$ws.Range("M22").Value = -179
# Row 31 (hunk 27)
$ws.Range("H31").Value = 44643.668
$ws.Range("I31").Value = 1614.7646
$ws.Range("K31").Value = 1614.7646
$ws.Range("M31").Value = -1319.7646
# Row 34 (hunk 28)
$ws.Range("H34").Value = 44643.668
$ws.Range("I34").Value = 1614.7646
$ws.Range("K34").Value = 1614.7646
$ws.Range("M34").Value = -1412.7646
# Row 86 (hunk 29)
$ws.Range("H86").Value = 3427.2964
$ws.Range("I86").Value = 3170.611
$ws.Range("J86").Value = 3940.6667
$ws.Range("K86").Value = 3170.611
$ws.Range("L86").Value = 3940.6667
$ws.Range("M86").Value = -2047.611
$ws.Range("N86").Value = -6186.6667
# Row 89 (hunk 30)
$ws.Range("H89").Value = 3427.2964
$ws.Range("I89").Value = 3170.611
$ws.Range("J89").Value = 3940.6667
$ws.Range("K89").Value = 15853.055
$ws.Range("L89").Value = 19703.3335
$ws.Range("M89").Value = -10237.055
$ws.Range("N89").Value = -30935.3335
# Row 113 (hunk 31)
$ws.Range("H113").Value = 3429.0417
$ws.Range("J113").Value = 3719.625
$ws.Range("L113").Value = 3719.625
$ws.Range("N113").Value = -8059.625
# Row 132 (hunk 32)
$ws.Range("H132").Value = 1695.2858
$ws.Range("I132").Value = 1491.3513
$ws.Range("K132").Value = 4474.0539
$ws.Range("M132").Value = -1944.0539
# Row 134 (hunk 33)
$ws.Range("H134").Value = 280614.75
$ws.Range("I134").Value = 2765.7812
$ws.Range("K134").Value = 8297.3436
$ws.Range("M134").Value = -5762.3436

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 56 (hunk 34)
$ws.Range("H56").Value = 7999.5
$ws.Range("I56").Value = 7999.5
$ws.Range("K56").Value = 7999.5
$ws.Range("M56").Value = -7469.5
# Row 80 (hunk 35)
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = ""
# Row 83 (hunk 36)
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = ""
# Row 101 (hunk 37)
$ws.Range("H101").Value = 5414
$ws.Range("J101").Value = 5914.5
$ws.Range("L101").Value = 17743.5
$ws.Range("N101").Value = -22611.5
# Row 113 (hunk 38)
$ws.Range("H113").Value = 1090757.1
$ws.Range("I113").Value = 3087512.2
$ws.Range("J113").Value = 1617.9546
$ws.Range("K113").Value = 9262536.600000001
$ws.Range("L113").Value = 4853.8638
$ws.Range("M113").Value = -9260366.600000001
$ws.Range("N113").Value = -9193.863799999999
# Row 119 (hunk 39)
$ws.Range("H119").Value = 10264.25
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 2 (hunk 40)
$ws.Range("H2").Value = 104.333336
$ws.Range("I2").Value = 51.333332
$ws.Range("J2").Value = 263.33334
$ws.Range("K2").Value = 51.333332
$ws.Range("L2").Value = 263.33334
$ws.Range("M2").Value = 61.666668
$ws.Range("N2").Value = -489.33334
# Row 132 (hunk 41)
$ws.Range("H132").Value = 88050.46000000001
$ws.Range("I132").Value = 13564.3
$ws.Range("J132").Value = 336337.66
$ws.Range("K132").Value = 40692.89999999999
$ws.Range("L132").Value = 1009012.98
$ws.Range("M132").Value = -38162.89999999999
$ws.Range("N132").Value = -1014072.98

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22 (hunk 42)
$ws.Range("H22").Value = 366.33334
$ws.Range("I22").Value = 299.5
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 299.5
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = -4.5
$ws.Range("N22").Value = -1090
# Row 23 (hunk 43)
$ws.Range("H23").Value = 732500
$ws.Range("I23").Value = 835714.3
$ws.Range("K23").Value = 835714.3
$ws.Range("M23").Value = -835484.3
# Row 27 (hunk 44)
$ws.Range("H27").Value = 366.33334
$ws.Range("I27").Value = 299.5
$ws.Range("J27").Value = 500
$ws.Range("K27").Value = 299.5
$ws.Range("L27").Value = 500
$ws.Range("M27").Value = -192.5
$ws.Range("N27").Value = -714
# Row 46 (hunk 45)
$ws.Range("H46").Value = 1998.3334
$ws.Range("J46").Value = 1738.4
$ws.Range("L46").Value = 1738.4
$ws.Range("N46").Value = -2114.4
# Row 136 (hunk 46)
$ws.Range("H136").Value = 320121.8
$ws.Range("I136").Value = 594140.0600000001
$ws.Range("J136").Value = 9567.799999999999
$ws.Range("K136").Value = 1782420.18
$ws.Range("L136").Value = 28703.4
$ws.Range("M136").Value = -1779870.18
$ws.Range("N136").Value = -33803.39999999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 96 (hunk 47)
$ws.Range("H96").Value = 112341.89
$ws.Range("I96").Value = 333999.66
$ws.Range("J96").Value = 1513
$ws.Range("K96").Value = 333999.66
$ws.Range("L96").Value = 1513
$ws.Range("M96").Value = -332626.66
$ws.Range("N96").Value = -4259
# Row 132 (hunk 48)
$ws.Range("H132").Value = 34711.56
$ws.Range("I132").Value = 5103.1924
$ws.Range("K132").Value = 15309.5772
$ws.Range("M132").Value = -12779.5772
# Row 136 (hunk 49)
$ws.Range("H136").Value = 92460.25999999999
$ws.Range("I136").Value = 22866.047
$ws.Range("J136").Value = 336040
$ws.Range("K136").Value = 68598.141
$ws.Range("L136").Value = 1008120
$ws.Range("M136").Value = -66048.141
$ws.Range("N136").Value = -1013220

